$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.610.80'
$ws.Range("E2").Value = '  +2.05%  '

$ws.Range("D3").Value = '1.888.50'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.30'
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4904'
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2951'
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06780'
$ws.Range("E9").Value = '  +2.17%  '

$ws.Range("D10").Value = '1.885.82'
$ws.Range("E10").Value = '  +0.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.23'
$ws.Range("E11").Value = '  +2.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07244'
$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '90.96'
$ws.Range("E13").Value = '  +4.17%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.061'
$ws.Range("E14").Value = '  +3.36%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6764'
$ws.Range("E15").Value = '  +1.16%  '

$ws.Range("D16").Value = '30.597.16'
$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007955'
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.13'
$ws.Range("E19").Value = '  +2.50%  '

$ws.Range("D20").Value = '2.130.95'
$ws.Range("E20").Value = '  +0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.825'
$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '187.64'
$ws.Range("E23").Value = '  +32.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.083'
$ws.Range("E24").Value = '  +3.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.340'
$ws.Range("E25").Value = '  +2.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.79'
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.10'
$ws.Range("E27").Value = '  +12.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.904'
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("E29").Value = '  +1.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.336'
$ws.Range("E30").Value = '  +2.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09048'
$ws.Range("E31").Value = '  +3.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.011'
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05194'
$ws.Range("E33").Value = '  +2.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7514'
$ws.Range("E34").Value = '  +4.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.752'
$ws.Range("E36").Value = '  +3.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01846'
$ws.Range("E37").Value = '  +3.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.664'
$ws.Range("E38").Value = '  -1.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.142'
$ws.Range("E39").Value = '  -1.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9346'
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4424'
$ws.Range("E41").Value = '  +4.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.17'
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.731'
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.583'
$ws.Range("E45").Value = '  +2.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1337'
$ws.Range("E46").Value = '  +5.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05851'
$ws.Range("E47").Value = '  +2.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.441'
$ws.Range("E48").Value = '  +7.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.699'
$ws.Range("E49").Value = '  +4.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3929'
$ws.Range("E50").Value = '  +4.09%  '

$ws.Range("E51").Value = '  +2.49%  '
